# Commit: "Optie toegevoegd om verschillen die na afronding gelijk lijken te verbergen."
# (Option added to hide differences that appear equal after rounding.)
#
# Two new configuration columns are introduced:
#   1. "indeling_rijen" sheet: new column "hoger_is_beter" inserted right after
#      "inhoud" (col C) - marked TRUE for the three "var" rows that previously
#      had no third column value (rows 10, 11, 12).
#   2. "algemeen" sheet: new column "sign_verbergen_wanneer_afgerond_gelijk"
#      inserted right after "sign_hovertekst" (col J) with default value FALSE.

$wb = $excel.ActiveWorkbook

# --- indeling_rijen: insert "hoger_is_beter" column C -----------------------
$wsRijen = $wb.Worksheets.Item("indeling_rijen")
$wsRijen.Columns.Item(3).Insert()
$wsRijen.Cells.Item(1, 3).Value = "hoger_is_beter"
$wsRijen.Cells.Item(10, 3).Value = $true
$wsRijen.Cells.Item(11, 3).Value = $true
$wsRijen.Cells.Item(12, 3).Value = $true
$wsRijen.Range("D19").Select()

# --- algemeen: insert "sign_verbergen_wanneer_afgerond_gelijk" column J ----
$wsAlgemeen = $wb.Worksheets.Item("algemeen")
$wsAlgemeen.Columns.Item(10).Insert()
$wsAlgemeen.Cells.Item(1, 10).Value = "sign_verbergen_wanneer_afgerond_gelijk"
$wsAlgemeen.Cells.Item(2, 10).Value = $false

# Make "algemeen" the active sheet/tab with J1 selected (last action wins).
$wsAlgemeen.Range("J1").Select()
